# Auto-generated Excel COM-interop script
# Updated cryptos list on Tue Jul 30 05:31:07 UTC 2024 with GitHub Actions

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.480.19"
$ws.Range("E2").Value = "  -4.41%  "
$ws.Range("D3").Value = "3.303.94"
$ws.Range("E3").Value = "  -1.29%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").Value = "'571.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.49%  "
$ws.Range("B6").Value = "Solana"
$ws.Range("C6").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D6").Value = "'182.03"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.77%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "'0.598"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.46%  "
$ws.Range("E9").Value = "  -3.83%  "
$ws.Range("D10").Value = "'6.62"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.17%  "
$ws.Range("D11").Value = "'0.402"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.70%  "
$ws.Range("D12").Value = "3.880.71"
$ws.Range("E12").Value = "  -1.17%  "
$ws.Range("E13").Value = "  -0.75%  "
$ws.Range("D14").Value = "'27.10"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.41%  "
$ws.Range("D15").Value = "66.566.04"
$ws.Range("E15").Value = "  -4.26%  "
$ws.Range("E16").Value = "  -2.85%  "
$ws.Range("D17").Value = "3.309.17"
$ws.Range("E17").Value = "  -1.20%  "
$ws.Range("D18").Value = "'13.65"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.45%  "
$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D19").Value = "'5.67"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.68%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "'429.73"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.06%  "
$ws.Range("D21").Value = "'7.61"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.05%  "
$ws.Range("D22").Value = "'73.58"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.16%  "
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("D24").Value = "'0.515"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.28%  "
$ws.Range("E25").Value = "  -3.46%  "
$ws.Range("E26").Value = "  +0.43%  "
$ws.Range("D27").Value = "'9.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.92%  "
$ws.Range("E28").Value = "  -0.88%  "
$ws.Range("E29").Value = "  -2.19%  "
$ws.Range("D30").Value = "'22.74"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.51%  "
$ws.Range("D31").Value = "'5.31"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.99%  "
$ws.Range("D32").Value = "'0.999"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.01%  "
$ws.Range("D33").Value = "'1.23"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.38%  "
$ws.Range("D34").Value = "'6.76"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.92%  "
$ws.Range("D36").Value = "'159.92"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.71%  "
$ws.Range("D37").Value = "'1.84"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Value = "'27.11"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.08%  "
$ws.Range("D39").Value = "2.792.10"
$ws.Range("E39").Value = "  +1.24%  "
$ws.Range("D40").Value = "'0.788"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.91%  "
$ws.Range("D41").Value = "'4.43"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.12%  "
$ws.Range("D42").Value = "'6.17"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.06%  "
$ws.Range("D43").Value = "'0.0673"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.24%  "
$ws.Range("D44").Value = "'40.16"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.28%  "
$ws.Range("D45").Value = "'24.31"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.99%  "
$ws.Range("E46").Value = "  -6.93%  "
$ws.Range("D47").Value = "'320.23"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.75%  "
$ws.Range("D48").Value = "'0.0271"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.22%  "
$ws.Range("D49").Value = "'0.981"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.76%  "
$ws.Range("E50").Value = "  -2.04%  "
$ws.Range("E51").Value = "  -1.21%  "
